$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns
# Note: Price values that look like plain numbers are prefixed with an
# apostrophe so Excel keeps storing them as text, matching the source data.
$ws.Range("D2").Value = "21.682.94"
$ws.Range("E2").Value = "  -1.56%  "
$ws.Range("D3").Value = "1.536.64"
$ws.Range("E3").Value = "  -1.08%  "
$ws.Range("D4").Value = "'0.9983"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "'0.9994"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").Value = "'289.90"
$ws.Range("E6").Value = "  +1.10%  "
$ws.Range("D7").Value = "'0.3941"
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("D8").Value = "'0.3174"
$ws.Range("E8").Value = "  -1.01%  "
$ws.Range("D9").Value = "'42.40"
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("D10").Value = "'0.07188"
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("D11").Value = "'1.076"
$ws.Range("E11").Value = "  -1.54%  "
$ws.Range("D12").Value = "'0.9983"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").Value = "'5.729"
$ws.Range("E13").Value = "  +1.74%  "
$ws.Range("D14").Value = "'18.50"
$ws.Range("E14").Value = "  -1.48%  "
$ws.Range("D15").Value = "'6.636"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").Value = "1.536.18"
$ws.Range("E16").Value = "  -0.94%  "
$ws.Range("D17").Value = "'0.00001097"
$ws.Range("E17").Value = "  -2.48%  "
$ws.Range("D18").Value = "'0.06598"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "'84.33"
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("D20").Value = "'0.9992"
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").Value = "'6.160"
$ws.Range("E21").Value = "  -2.49%  "
$ws.Range("D22").Value = "'15.59"
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("D23").Value = "'10.76"
$ws.Range("E23").Value = "  -3.93%  "
$ws.Range("D24").Value = "'2.375"
$ws.Range("E24").Value = "  +1.19%  "
$ws.Range("D25").Value = "21.683.48"
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("D26").Value = "'2.380"
$ws.Range("E26").Value = "  -2.67%  "
$ws.Range("D27").Value = "'151.22"
$ws.Range("E27").Value = "  +2.54%  "
$ws.Range("D28").Value = "'18.44"
$ws.Range("E28").Value = "  -1.20%  "
$ws.Range("D29").Value = "'4.865"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").Value = "1.709.17"
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("D31").Value = "'117.55"
$ws.Range("E31").Value = "  -1.98%  "
$ws.Range("D32").Value = "'6.127"
$ws.Range("E32").Value = "  +8.26%  "
$ws.Range("D33").Value = "'0.9634"
$ws.Range("E33").Value = "  -8.98%  "
$ws.Range("D35").Value = "'5.219"
$ws.Range("E35").Value = "  +2.39%  "
$ws.Range("D36").Value = "'8.549"
$ws.Range("E36").Value = "  -7.34%  "
$ws.Range("D39").Value = "'1.466"
$ws.Range("E39").Value = "  -8.15%  "
$ws.Range("D42").Value = "'1.187"
$ws.Range("E42").Value = "  -1.40%  "
$ws.Range("D43").Value = "'0.9988"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("D44").Value = "'0.5848"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").Value = "'13.10"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("D46").Value = "'3.730"
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("D47").Value = "'0.5589"
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("D48").Value = "'1.170"
$ws.Range("E48").Value = "  +2.53%  "
$ws.Range("D49").Value = "'1.891"
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("D50").Value = "'116.85"
$ws.Range("E50").Value = "  -0.75%  "
$ws.Range("D51").Value = "'0.06710"
$ws.Range("E51").Value = "  -1.89%  "

# Row 34: only Volume(1h) changes
$ws.Range("E34").Value = "  -2.19%  "

# Rows 37/38 swap: Hedera <-> VeChain (with updated price/volume)
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02236"
$ws.Range("E37").Value = "  -1.16%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.06034"
$ws.Range("E38").Value = "  -2.49%  "

# Rows 40/41 swap: Algorand <-> Aptos (with updated price/volume)
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "'11.20"
$ws.Range("E40").Value = "  +5.91%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "'0.2043"
$ws.Range("E41").Value = "  -1.96%  "
